# Update crypto price/volume snapshot cells (Price = column D, Volume(1h) = column E).
# Cells are plain text in the source sheet; for D-column values that look like a
# plain number, force text via NumberFormat "@" then restore the original "Normal"
# cell style so no visible formatting/style index changes leak into the saved file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.842.03"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.084.09"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.82%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0789"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").Value = "2.392.26"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.774"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "2.098.57"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "37.773.52"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "0.0₃0850"
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0634"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.45%  "
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").Value = "1.449.14"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").Value = "2.277.32"
